$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 29, shifting existing rows 29-42 down to 31-44
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(29).Insert()

# Row 29: Angeleno / Primera
$ws.Cells.Item(29,1).Value = 2
$ws.Cells.Item(29,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(29,3).Value = 'Coquimbo'
$ws.Cells.Item(29,4).Value = 44651
$ws.Cells.Item(29,5).Value = 4
$ws.Cells.Item(29,6).Value = 'Fruta'
$ws.Cells.Item(29,7).Value = 100103
$ws.Cells.Item(29,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(29,9).Value = 100103002
$ws.Cells.Item(29,10).Value = 'Ciruela'
$ws.Cells.Item(29,11).Value = 'Angeleno'
$ws.Cells.Item(29,12).Value = 'Primera'
$ws.Cells.Item(29,13).Value = 10
$ws.Cells.Item(29,14).Value = 210000
$ws.Cells.Item(29,15).Value = 220000
$ws.Cells.Item(29,16).Value = 215000
$ws.Cells.Item(29,17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(29,18).Value = "Región de O'Higgins"
$ws.Cells.Item(29,19).Value = 478
$ws.Cells.Item(29,20).Value = 450

# Row 30: Angeleno / Segunda
$ws.Cells.Item(30,1).Value = 2
$ws.Cells.Item(30,2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(30,3).Value = 'Coquimbo'
$ws.Cells.Item(30,4).Value = 44651
$ws.Cells.Item(30,5).Value = 4
$ws.Cells.Item(30,6).Value = 'Fruta'
$ws.Cells.Item(30,7).Value = 100103
$ws.Cells.Item(30,8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(30,9).Value = 100103002
$ws.Cells.Item(30,10).Value = 'Ciruela'
$ws.Cells.Item(30,11).Value = 'Angeleno'
$ws.Cells.Item(30,12).Value = 'Segunda'
$ws.Cells.Item(30,13).Value = 16
$ws.Cells.Item(30,14).Value = 170000
$ws.Cells.Item(30,15).Value = 180000
$ws.Cells.Item(30,16).Value = 175000
$ws.Cells.Item(30,17).Value = '$/bins (450 kilos)'
$ws.Cells.Item(30,18).Value = "Región de O'Higgins"
$ws.Cells.Item(30,19).Value = 389
$ws.Cells.Item(30,20).Value = 450

Write-Output "Done"
